$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "pair_kind" column (J) for the practice rows 2-5 with "generic"
$ws.Cells.Item(2, 10).Value = "generic"
$ws.Cells.Item(3, 10).Value = "generic"
$ws.Cells.Item(4, 10).Value = "generic"
$ws.Cells.Item(5, 10).Value = "generic"

# New "stim details" section starting at row 27
$ws.Cells.Item(27, 1).Value = "stim details"

# Header row for the stim details table
$ws.Cells.Item(28, 1).Value = "month"
$ws.Cells.Item(28, 2).Value = "word_type"
$ws.Cells.Item(28, 3).Value = "need_audio"
$ws.Cells.Item(28, 4).Value = "need_image"
$ws.Cells.Item(28, 5).Value = "word"
$ws.Cells.Item(28, 6).Value = "count"
$ws.Cells.Item(28, 7).Value = "find images"

# Data rows 29-36
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "video"

$ws.Cells.Item(30, 1).Value = 6
$ws.Cells.Item(30, 2).Value = "video"

$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "video"

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "video"

$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "audio"

$ws.Cells.Item(34, 1).Value = 6
$ws.Cells.Item(34, 2).Value = "audio"

$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "audio"

$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "audio"
